# Update "想去人数" (want-to-go count) figures in column F across the
# workbook's sheets to match a refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheetId 1) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 193
$ws.Range("F4").Value = 193
$ws.Range("F5").Value = 46
$ws.Range("F6").Value = 1603
$ws.Range("F7").Value = 3212
$ws.Range("F8").Value = 737
$ws.Range("F9").Value = 1954
$ws.Range("F10").Value = 1878
$ws.Range("F11").Value = 957
$ws.Range("F12").Value = 333
$ws.Range("F13").Value = 12
$ws.Range("F14").Value = 1563
$ws.Range("F15").Value = 327
$ws.Range("F17").Value = 54
$ws.Range("F18").Value = 1366
$ws.Range("F19").Value = 463
$ws.Range("F20").Value = 573
$ws.Range("F21").Value = 268
$ws.Range("F22").Value = 10293
$ws.Range("F23").Value = 9490
$ws.Range("F24").Value = 819
$ws.Range("F25").Value = 627
$ws.Range("F26").Value = 1791
$ws.Range("F27").Value = 125
$ws.Range("F28").Value = 348

# --- Sheet "演出" (sheetId 2) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 34
$ws.Range("F5").Value = 110

# --- Sheet "全部类型" (sheetId 4) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 193
$ws.Range("F5").Value = 193
$ws.Range("F7").Value = 46
$ws.Range("F8").Value = 1603
$ws.Range("F9").Value = 3212
$ws.Range("F10").Value = 737
$ws.Range("F11").Value = 1954
$ws.Range("F12").Value = 1878
$ws.Range("F13").Value = 958
$ws.Range("F14").Value = 333
$ws.Range("F15").Value = 12
$ws.Range("F16").Value = 1563
$ws.Range("F17").Value = 327
$ws.Range("F20").Value = 54
$ws.Range("F21").Value = 34
$ws.Range("F22").Value = 1366
$ws.Range("F23").Value = 463
$ws.Range("F24").Value = 573
$ws.Range("F25").Value = 268
$ws.Range("F26").Value = 10293
$ws.Range("F27").Value = 9490
$ws.Range("F28").Value = 819
$ws.Range("F29").Value = 627
$ws.Range("F30").Value = 1791
$ws.Range("F31").Value = 110
$ws.Range("F33").Value = 125
$ws.Range("F34").Value = 348
